# #5: property building done
#
# Adds a data row (row 2) to every property-type worksheet (land, building,
# car, cash) and, for the deposit sheet, inserts a new row 2 ahead of the
# two rows that were already there. The building sheet's old single data
# row (which used an ad-hoc column layout) is replaced by the common
# 16-column header used by every other sheet, and its data is moved into
# row 2 under that layout.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 土地 land
$ws2 = $wb.Worksheets.Item(2)  # 建物 building
$ws3 = $wb.Worksheets.Item(3)  # 汽車 car
$ws4 = $wb.Worksheets.Item(4)  # 現金 cash
$ws5 = $wb.Worksheets.Item(5)  # 存款 deposit

# =========================================================================
# Sheet1 (土地 / land): add row 2
# =========================================================================

# Give A2 the same bold/bordered "index" look as the header row (style
# used by every other index cell in the workbook) by copying the format
# from the header before writing the value.
$ws1.Range("B1").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Range("A2").Value = 14

$ws1.Range("B2").Value = "新北市永和區永福段08230000地號"
$ws1.Range("C2").Value = 20
$ws1.Range("D2").Value = "5分之1"
$ws1.Range("E2").Value = "邱若山"
$ws1.Range("F2").Value = "80年06月13曰"
$ws1.Range("G2").Value = "第一次登記"
$ws1.Range("H2").Value = "(超過五年）"
$ws1.Range("I2").Value = "land"
$ws1.Range("J2").Value = "normal"
# "2013-12-30" looks like a date, write it then strip the auto-applied
# date formatting so it stays a plain text value.
$ws1.Range("K2").Value = "2013-12-30"
$ws1.Range("B2").Copy()
$ws1.Range("K2").PasteSpecial(-4122)
$ws1.Range("K2").Value = "2013-12-30"
$ws1.Range("L2").Value = "林淑芬"
$ws1.Range("M2").Value = 1337
$ws1.Range("N2").Value = "tmp63cf1"
$ws1.Range("O2").Value = 14
$ws1.Range("P2").Value = 0.2
$ws1.Range("Q2").Value = 4

# =========================================================================
# Sheet2 (建物 / building): replace row1 with the common header and move
# the old row1 data into row2 under the new layout.
# =========================================================================
$ws2.Range("B1").Value = "name"
$ws2.Range("C1").Value = "area"
$ws2.Range("D1").Value = "share_portion"
$ws2.Range("E1").Value = "owner"
$ws2.Range("F1").Value = "register_date"
$ws2.Range("G1").Value = "register_reason"
$ws2.Range("H1").Value = "acquire_value"
$ws2.Range("I1").Value = "property_category"
$ws2.Range("J1").Value = "category"
$ws2.Range("K1").Value = "date"
$ws2.Range("L1").Value = "legislator_name"
$ws2.Range("M1").Value = "legislator_id"
$ws2.Range("N1").Value = "source_file"
$ws2.Range("O1").Value = "index"
$ws2.Range("P1").Value = "portion"
$ws2.Range("Q1").Value = "total"

$ws2.Range("B1").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$ws2.Range("A2").Value = 19

$ws2.Range("B2").Value = "新北市永和區永福段02017000建號"
$ws2.Range("C2").Value = 113
$ws2.Range("D2").Value = "全部"
$ws2.Range("E2").Value = "邱若山"
$ws2.Range("F2").Value = "80年06月13日"
$ws2.Range("G2").Value = "第一次登記"
$ws2.Range("H2").Value = "(超過五年）"
$ws2.Range("I2").Value = "land"
$ws2.Range("J2").Value = "normal"
$ws2.Range("K2").Value = "2013-12-30"
$ws2.Range("B2").Copy()
$ws2.Range("K2").PasteSpecial(-4122)
$ws2.Range("K2").Value = "2013-12-30"
$ws2.Range("L2").Value = "林淑芬"
$ws2.Range("M2").Value = 1337
$ws2.Range("N2").Value = "tmp63cf1"
$ws2.Range("O2").Value = 19
$ws2.Range("P2").Value = 1
$ws2.Range("Q2").Value = 113

# =========================================================================
# Sheet3 (汽車 / car): add row2 duplicating row1 plus an index column.
# =========================================================================
$ws3.Range("B1").Copy()
$ws3.Range("A2").PasteSpecial(-4122)
$ws3.Range("A2").Value = 29

$ws3.Range("B2").Value = "toyotarav4rod"
$ws3.Range("C2").Value = 2362
$ws3.Range("D2").Value = "林淑芬"
$ws3.Range("E2").Value = "100年03月01曰"
$ws3.Range("F2").Value = "買賣"
$ws3.Range("G2").Value = 989000

# =========================================================================
# Sheet4 (現金 / cash): add row2 duplicating row1 plus an index column.
# =========================================================================
$ws4.Range("B1").Copy()
$ws4.Range("A2").PasteSpecial(-4122)
$ws4.Range("A2").Value = 41

$ws4.Range("B2").Value = "新臺幣"
$ws4.Range("C2").Value = "林淑芬"
$ws4.Range("D2").Value = 2320000

# =========================================================================
# Sheet5 (存款 / deposit): insert a new row2 (shifts the old row2/row3
# down to row3/row4) duplicating the data already on row1.
# =========================================================================
$ws5.Rows(2).Insert()

$ws5.Range("A3").Copy()
$ws5.Range("A2").PasteSpecial(-4122)
$ws5.Range("A2").Value = 46

$ws5.Range("B3:F3").Copy()
$ws5.Range("B2:F2").PasteSpecial(-4122)

$ws5.Range("B2").Value = "臺灣銀行群賢分行"
$ws5.Range("C2").Value = "活期存款"
$ws5.Range("D2").Value = "新臺幣"
$ws5.Range("E2").Value = "林淑芬"
$ws5.Range("F2").Value = 900000
